# Applies the 2025-11-14 06:27 JST scrape-append update to the 案件情報 workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Remove the 7 oldest listings (former rows 10-16); remaining rows shift up ---
$ws.Rows("10:16").Delete()

# --- Refresh the 8 surviving rows (2-9) with the newer scrape snapshot ---
# row 2
$ws.Range("A2").Value = '2025-11-14 06:27:31'
$ws.Range("B2").Value = '【GAS】Yahoo!ショッピング注文完了メール (Gmail) からスプレッドシートに転記する仕事'
$ws.Range("D2").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("G2").Value = 298
$ws.Range("H2").Value = '🔥AI,Ai'

# row 3
$ws.Range("A3").Value = '2025-11-14 06:27:31'
$ws.Range("B3").Value = '英語教育の公式LINEアカウント開発・運用スタッフ募集【即日〜3月/4ヶ月/継続可能】'
$ws.Range("D3").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("G3").Value = 63
$ws.Range("H3").Value = '◆開発'

# row 4
$ws.Range("A4").Value = '2025-11-14 06:27:31'
$ws.Range("B4").Value = '【急募】WordPressでの商品検索サイト構築依頼'
$ws.Range("D4").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("G4").Value = 58
$ws.Range("H4").Value = '◇サイト ○WordPress'

# row 5
$ws.Range("A5").Value = '2025-11-14 06:27:31'
$ws.Range("B5").Value = 'wordpressレンダリングを妨げるリソースの除外'
$ws.Range("D5").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("G5").Value = 33
$ws.Range("H5").Value = '○WordPress'

# row 6
$ws.Range("A6").Value = '2025-11-14 06:27:31'
$ws.Range("B6").Value = '【相談から】Laravel7からLaravel12へのサーバーアップデート依頼'
$ws.Range("D6").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("G6").Value = 25
$ws.Range("H6").ClearContents()

# row 7
$ws.Range("A7").Value = '2025-11-14 06:27:31'
$ws.Range("B7").Value = '初回 ★社内の音響設計スキル向上のため、Modeler / EASE Focus を教えていただける方'
$ws.Range("D7").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("G7").Value = 18
$ws.Range("H7").ClearContents()

# row 8
$ws.Range("A8").Value = '2025-11-14 06:27:31'
$ws.Range("B8").Value = '月1~5万円以内の小規模タスク依頼'
$ws.Range("D8").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("G8").Value = 13
$ws.Range("H8").ClearContents()

# row 9
$ws.Range("A9").Value = '2025-11-14 06:27:31'
$ws.Range("B9").Value = '【急募】Unityで自動ルート設計プログラムが実現可能か専門家に相談がしたい'
$ws.Range("D9").Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range("G9").Value = 10
$ws.Range("H9").ClearContents()

# --- Rebuild hyperlinks on column F cleanly (avoids orphaned/duplicate links) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5433649') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5433668') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5433985') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5016989') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5433727') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5433823') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5433937') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5434061') | Out-Null

# --- Column width tweaks (Excel stores ColumnWidth + 5/6 char as the OOXML width) ---
$colWidthPad = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 53 - $colWidthPad   # B: 52 -> 53
$ws.Columns.Item(4).ColumnWidth = 30 - $colWidthPad   # D: 32 -> 30
$ws.Columns.Item(8).ColumnWidth = 17 - $colWidthPad   # H: 27 -> 17

